# Add a new "true" (text) column C to the existing testdata rows,
# entered with a leading apostrophe so Excel stores it as quote-prefixed text
# rather than as a boolean value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 5; $row++) {
    $ws.Cells.Item($row, 3).Value = "'true"
}

# Column A grew slightly wider to accommodate the new layout.
$ws.Columns.Item(1).ColumnWidth = 22.5
